$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 4.6
$ws.Range("J2").Value = 4.1
$ws.Range("K2").Value = 4.3
$ws.Range("L2").Value = 1.34
$ws.Range("V2").Value = 1.27
$ws.Range("X2").Value = 19.5
$ws.Range("Z2").Value = 40
$ws.Range("AA2").Value = 100
$ws.Range("AC2").Value = 9.4
$ws.Range("AD2").Value = 18
$ws.Range("AE2").Value = 55
$ws.Range("AF2").Value = 12.5
$ws.Range("AG2").Value = 10
$ws.Range("AH2").Value = 17.5
$ws.Range("AI2").Value = 55
$ws.Range("AK2").Value = 18
$ws.Range("AL2").Value = 30
$ws.Range("AM2").Value = 95
$ws.Range("AO2").Value = 46

# Row 3
$ws.Range("I3").Value = 2.18

# Row 4
$ws.Range("P4").Value = 2.02

# Row 5
$ws.Range("G5").Value = 4.9
$ws.Range("H5").Value = 1.7
$ws.Range("K5").Value = 5

# Row 6
$ws.Range("F6").Value = 2.92
$ws.Range("G6").Value = 2.96
$ws.Range("I6").Value = 2.7
$ws.Range("T6").Value = 1.78

# Row 7
$ws.Range("I7").Value = 1.58
$ws.Range("K7").Value = 4.5

# Row 8
$ws.Range("G8").Value = 2.3
$ws.Range("I8").Value = 4.4
$ws.Range("Y8").Value = 9.4
$ws.Range("AH8").Value = 38
